$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value2 = 45183
$ws.Range("C3").Value2 = 45183
$ws.Range("C4").Value2 = 45183
$ws.Range("C5").Value2 = 45183
